{"js": "// Replace equation text in each table cell according to the mapping below.\n// All \"old\" values are unique within the document (verified against the diff),\n// so a direct search+replace per pair is safe and order-independent.\nconst replacements = {\n  \"80-10=\": \"67-47=\",\n  \"45-10=\": \"89-71=\",\n  \"62-31=\": \"22+11=\",\n  \"40+21=\": \"29+16=\",\n  \"35+10=\": \"39+23=\",\n  \"75-24=\": \"86-1=\",\n  \"42+18=\": \"92-46=\",\n  \"56+12=\": \"83-34=\",\n  \"34+17=\": \"20+32=\",\n  \"97-10=\": \"36-1=\",\n  \"47-21=\": \"8+74=\",\n  \"76+22=\": \"9+65=\",\n  \"69+15=\": \"38+26=\",\n  \"35+63=\": \"16+61=\",\n  \"44+47=\": \"57+35=\",\n  \"79-68=\": \"42-11=\",\n  \"97-42=\": \"66+0=\",\n  \"5+70=\": \"76+1=\",\n  \"7+51=\": \"86-85=\",\n  \"23+76=\": \"6+91=\",\n  \"29+57=\": \"5+86=\",\n  \"7-4=\": \"69-54=\",\n  \"49+47=\": \"98-13=\",\n  \"69-34=\": \"2+13=\",\n  \"10+36=\": \"86+13=\",\n  \"1+95=\": \"20+76=\",\n  \"8+28=\": \"76-6=\",\n  \"78-37=\": \"4+75=\",\n  \"91+4=\": \"97-16=\",\n  \"33+33=\": \"82-3=\",\n  \"13+85=\": \"10+2=\",\n  \"77-52=\": \"87-74=\",\n  \"56+20=\": \"96-1=\",\n  \"32+47=\": \"49+5=\",\n  \"89-21=\": \"14-10=\",\n  \"68+10=\": \"4+50=\",\n  \"60+19=\": \"72-62=\",\n  \"31-23=\": \"93-57=\",\n  \"82-24=\": \"20+37=\",\n  \"40-35=\": \"74-57=\",\n  \"50-16=\": \"89-78=\",\n  \"8+77=\": \"94-71=\",\n  \"95-46=\": \"35-4=\",\n  \"55+27=\": \"82+13=\",\n  \"41+11=\": \"69-11=\",\n  \"87-67=\": \"60-48=\",\n  \"17+30=\": \"68-67=\",\n  \"75-1=\": \"31+45=\",\n  \"42+44=\": \"77-25=\",\n  \"11+15=\": \"2+58=\",\n  \"83+9=\": \"34+61=\",\n  \"69-52=\": \"59+20=\",\n  \"82-70=\": \"31-15=\",\n  \"89-33=\": \"1+47=\",\n  \"89-74=\": \"64-43=\",\n  \"26-22=\": \"68-68=\",\n  \"62+1=\": \"11-4=\",\n  \"71-9=\": \"27+63=\",\n  \"33+20=\": \"53-6=\",\n  \"72-23=\": \"76-16=\",\n  \"7+90=\": \"78+6=\",\n  \"21+60=\": \"81-4=\",\n  \"34-25=\": \"86-9=\",\n  \"9+17=\": \"41+38=\",\n  \"36+37=\": \"66-23=\",\n  \"25+27=\": \"67-64=\",\n  \"22-8=\": \"40+58=\",\n  \"55-27=\": \"38-2=\",\n  \"17+53=\": \"1+94=\",\n  \"10+54=\": \"4+54=\",\n  \"37-29=\": \"37-23=\",\n  \"63+17=\": \"87-85=\",\n  \"36-8=\": \"89-9=\",\n  \"27-14=\": \"94+0=\",\n  \"21+8=\": \"39+42=\",\n  \"25+63=\": \"13+71=\",\n  \"13+23=\": \"93-92=\",\n  \"61-23=\": \"74-52=\",\n  \"72-3=\": \"40+40=\",\n  \"30+36=\": \"98-48=\",\n  \"40+6=\": \"99-89=\",\n  \"65-4=\": \"9+48=\",\n  \"55-45=\": \"57+15=\",\n  \"59-47=\": \"71-64=\",\n  \"11+67=\": \"71-62=\",\n  \"27+42=\": \"39+52=\",\n  \"92-33=\": \"91-32=\",\n  \"86-57=\": \"98-11=\",\n  \"30+29=\": \"93-50=\",\n  \"62-62=\": \"43-5=\",\n  \"28+65=\": \"60-34=\",\n  \"26+56=\": \"7+52=\",\n  \"38+42=\": \"46+10=\",\n  \"92-51=\": \"30-12=\",\n  \"10+74=\": \"66+17=\",\n  \"47-2=\": \"16+53=\",\n  \"64+15=\": \"31-1=\",\n  \"55-11=\": \"15+12=\",\n  \"67-17=\": \"34+38=\",\n  \"83-72=\": \"35+3=\"\n};\n\nconst body = context.document.body;\n\nfor (const oldText of Object.keys(replacements)) {\n  const newText = replacements[oldText];\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each equation's text in the table according to the mapping below.\n# All \"old\" values are unique within the document, so Find/Replace (wdReplaceAll,\n# scoped to each exact phrase) safely retargets exactly one run each, in any order.\n$d = $word.ActiveDocument\n\n$pairs = [ordered]@{\n    \"80-10=\" = \"67-47=\"\n    \"45-10=\" = \"89-71=\"\n    \"62-31=\" = \"22+11=\"\n    \"40+21=\" = \"29+16=\"\n    \"35+10=\" = \"39+23=\"\n    \"75-24=\" = \"86-1=\"\n    \"42+18=\" = \"92-46=\"\n    \"56+12=\" = \"83-34=\"\n    \"34+17=\" = \"20+32=\"\n    \"97-10=\" = \"36-1=\"\n    \"47-21=\" = \"8+74=\"\n    \"76+22=\" = \"9+65=\"\n    \"69+15=\" = \"38+26=\"\n    \"35+63=\" = \"16+61=\"\n    \"44+47=\" = \"57+35=\"\n    \"79-68=\" = \"42-11=\"\n    \"97-42=\" = \"66+0=\"\n    \"5+70=\" = \"76+1=\"\n    \"7+51=\" = \"86-85=\"\n    \"23+76=\" = \"6+91=\"\n    \"29+57=\" = \"5+86=\"\n    \"7-4=\" = \"69-54=\"\n    \"49+47=\" = \"98-13=\"\n    \"69-34=\" = \"2+13=\"\n    \"10+36=\" = \"86+13=\"\n    \"1+95=\" = \"20+76=\"\n    \"8+28=\" = \"76-6=\"\n    \"78-37=\" = \"4+75=\"\n    \"91+4=\" = \"97-16=\"\n    \"33+33=\" = \"82-3=\"\n    \"13+85=\" = \"10+2=\"\n    \"77-52=\" = \"87-74=\"\n    \"56+20=\" = \"96-1=\"\n    \"32+47=\" = \"49+5=\"\n    \"89-21=\" = \"14-10=\"\n    \"68+10=\" = \"4+50=\"\n    \"60+19=\" = \"72-62=\"\n    \"31-23=\" = \"93-57=\"\n    \"82-24=\" = \"20+37=\"\n    \"40-35=\" = \"74-57=\"\n    \"50-16=\" = \"89-78=\"\n    \"8+77=\" = \"94-71=\"\n    \"95-46=\" = \"35-4=\"\n    \"55+27=\" = \"82+13=\"\n    \"41+11=\" = \"69-11=\"\n    \"87-67=\" = \"60-48=\"\n    \"17+30=\" = \"68-67=\"\n    \"75-1=\" = \"31+45=\"\n    \"42+44=\" = \"77-25=\"\n    \"11+15=\" = \"2+58=\"\n    \"83+9=\" = \"34+61=\"\n    \"69-52=\" = \"59+20=\"\n    \"82-70=\" = \"31-15=\"\n    \"89-33=\" = \"1+47=\"\n    \"89-74=\" = \"64-43=\"\n    \"26-22=\" = \"68-68=\"\n    \"62+1=\" = \"11-4=\"\n    \"71-9=\" = \"27+63=\"\n    \"33+20=\" = \"53-6=\"\n    \"72-23=\" = \"76-16=\"\n    \"7+90=\" = \"78+6=\"\n    \"21+60=\" = \"81-4=\"\n    \"34-25=\" = \"86-9=\"\n    \"9+17=\" = \"41+38=\"\n    \"36+37=\" = \"66-23=\"\n    \"25+27=\" = \"67-64=\"\n    \"22-8=\" = \"40+58=\"\n    \"55-27=\" = \"38-2=\"\n    \"17+53=\" = \"1+94=\"\n    \"10+54=\" = \"4+54=\"\n    \"37-29=\" = \"37-23=\"\n    \"63+17=\" = \"87-85=\"\n    \"36-8=\" = \"89-9=\"\n    \"27-14=\" = \"94+0=\"\n    \"21+8=\" = \"39+42=\"\n    \"25+63=\" = \"13+71=\"\n    \"13+23=\" = \"93-92=\"\n    \"61-23=\" = \"74-52=\"\n    \"72-3=\" = \"40+40=\"\n    \"30+36=\" = \"98-48=\"\n    \"40+6=\" = \"99-89=\"\n    \"65-4=\" = \"9+48=\"\n    \"55-45=\" = \"57+15=\"\n    \"59-47=\" = \"71-64=\"\n    \"11+67=\" = \"71-62=\"\n    \"27+42=\" = \"39+52=\"\n    \"92-33=\" = \"91-32=\"\n    \"86-57=\" = \"98-11=\"\n    \"30+29=\" = \"93-50=\"\n    \"62-62=\" = \"43-5=\"\n    \"28+65=\" = \"60-34=\"\n    \"26+56=\" = \"7+52=\"\n    \"38+42=\" = \"46+10=\"\n    \"92-51=\" = \"30-12=\"\n    \"10+74=\" = \"66+17=\"\n    \"47-2=\" = \"16+53=\"\n    \"64+15=\" = \"31-1=\"\n    \"55-11=\" = \"15+12=\"\n    \"67-17=\" = \"34+38=\"\n    \"83-72=\" = \"35+3=\"\n}\n\nforeach ($old in $pairs.Keys) {\n    $new = $pairs[$old]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"No match found for '$old'\"\n    }\n}\n"}
